$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 11,15
$data[0,0] = 'TOTAL (5-17 ans)'
$data[0,1] = 'Tous les groupes de population'
$data[0,2] = 2307204
$data[0,3] = 21.3
$data[0,4] = 491129
$data[0,5] = 2.2
$data[0,6] = 49876
$data[0,7] = 26.2
$data[0,8] = 603526
$data[0,9] = 2.9
$data[0,10] = 67720
$data[0,11] = 47.5
$data[0,12] = 1094955
$data[0,13] = 52.5
$data[0,14] = 1212250
$data[1,0] = 'hote (5-17 y.o.)'
$data[1,1] = 'hote'
$data[1,2] = 2077436
$data[1,3] = 20.9
$data[1,4] = 434615
$data[1,5] = 2.1
$data[1,6] = 43979
$data[1,7] = 26
$data[1,8] = 540693
$data[1,9] = 2.9
$data[1,10] = 59824
$data[1,11] = 48.1
$data[1,12] = 998326
$data[1,13] = 51.9
$data[1,14] = 1079110
$data[2,0] = 'idp_host (5-17 y.o.)'
$data[2,1] = 'idp_host'
$data[2,2] = 125059
$data[2,3] = 23.8
$data[2,4] = 29787
$data[2,5] = 2
$data[2,6] = 2561
$data[2,7] = 24.7
$data[2,8] = 30925
$data[2,9] = 2.8
$data[2,10] = 3476
$data[2,11] = 46.6
$data[2,12] = 58310
$data[2,13] = 53.4
$data[2,14] = 66749
$data[3,0] = 'retourne (5-17 y.o.)'
$data[3,1] = 'retourne'
$data[3,2] = 77546
$data[3,3] = 24.3
$data[3,4] = 18874
$data[3,5] = 3.4
$data[3,6] = 2658
$data[3,7] = 31.4
$data[3,8] = 24380
$data[3,9] = 4.5
$data[3,10] = 3509
$data[3,11] = 36.3
$data[3,12] = 28124
$data[3,13] = 63.7
$data[3,14] = 49422
$data[4,0] = 'idp_site (5-17 y.o.)'
$data[4,1] = 'idp_site'
$data[4,2] = 27164
$data[4,3] = 28.9
$data[4,4] = 7853
$data[4,5] = 2.5
$data[4,6] = 678
$data[4,7] = 27.7
$data[4,8] = 7528
$data[4,9] = 3.4
$data[4,10] = 911
$data[4,11] = 37.5
$data[4,12] = 10195
$data[4,13] = 62.5
$data[4,14] = 16970
$data[5,0] = 'Filles (5-17 ans)'
$data[5,1] = 'Tous les groupes de population'
$data[5,2] = 1176674
$data[5,3] = 22.4
$data[5,4] = 263952
$data[5,5] = 2.4
$data[5,6] = 28789
$data[5,7] = 26
$data[5,8] = 306362
$data[5,9] = 2.9
$data[5,10] = 34418
$data[5,11] = 46.2
$data[5,12] = 543153
$data[5,13] = 53.8
$data[5,14] = 633521
$data[6,0] = 'Garcons (5-17 ans)'
$data[6,1] = 'Tous les groupes de population'
$data[6,2] = 1129288
$data[6,3] = 20.2
$data[6,4] = 228254
$data[6,5] = 1.9
$data[6,6] = 21528
$data[6,7] = 26.2
$data[6,8] = 295484
$data[6,9] = 2.9
$data[6,10] = 33235
$data[6,11] = 48.8
$data[6,12] = 550787
$data[6,13] = 51.2
$data[6,14] = 578500
$data[7,0] = 'Éducation préscolaire (5 ans)'
$data[7,1] = 'Tous les groupes de population'
$data[7,2] = 134153
$data[7,3] = 53.3
$data[7,4] = 71504
$data[7,5] = 3.7
$data[7,6] = 4967
$data[7,7] = 11.9
$data[7,8] = 15922
$data[7,9] = 1.5
$data[7,10] = 1948
$data[7,11] = 29.7
$data[7,12] = 39813
$data[7,13] = 70.3
$data[7,14] = 94340
$data[8,0] = 'École primaire'
$data[8,1] = 'Tous les groupes de population'
$data[8,2] = 1242341
$data[8,3] = 20.3
$data[8,4] = 252799
$data[8,5] = 1.5
$data[8,6] = 18596
$data[8,7] = 26.7
$data[8,8] = 331596
$data[8,9] = 3
$data[8,10] = 36994
$data[8,11] = 48.5
$data[8,12] = 602355
$data[8,13] = 51.5
$data[8,14] = 639986
$data[9,0] = 'Niveau scolaire intermédiaire'
$data[9,1] = 'Tous les groupes de population'
$data[9,2] = 709129
$data[9,3] = 17.4
$data[9,4] = 123654
$data[9,5] = 2.6
$data[9,6] = 18254
$data[9,7] = 26.9
$data[9,8] = 190763
$data[9,9] = 2.9
$data[9,10] = 20858
$data[9,11] = 50.1
$data[9,12] = 355600
$data[9,13] = 49.9
$data[9,14] = 353529
$data[10,0] = 'École secondaire'
$data[10,1] = 'Tous les groupes de population'
$data[10,2] = 175507
$data[10,3] = 19.4
$data[10,4] = 34056
$data[10,5] = 5.6
$data[10,6] = 9890
$data[10,7] = 26.2
$data[10,8] = 45991
$data[10,9] = 4.5
$data[10,10] = 7850
$data[10,11] = 44.3
$data[10,12] = 77721
$data[10,13] = 55.7
$data[10,14] = 97786

$ws.Range("A2:O12").Value = $data

